$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value2 = 0.2768050460063307
$ws.Range("C2").Value2 = 0.06031819443467157
$ws.Range("D2").Value2 = 0.07874435808558644
$ws.Range("E2").Value2 = 0.1686307853009623
$ws.Range("G2").Value2 = 0.5632483727616489
$ws.Range("H2").Value2 = 0.6924185566122532
$ws.Range("K2").Value2 = 0.2606336659522697
$ws.Range("M2").Value2 = 0.2148598588092483
$ws.Range("O2").Value2 = 2.48644661295495
$ws.Range("B3").Value2 = 0.2438897404451836
$ws.Range("C3").Value2 = 0.05720082609634858
$ws.Range("D3").Value2 = 0.07140343440914876
$ws.Range("E3").Value2 = 0.1576750876642521
$ws.Range("G3").Value2 = 0.5631517917375035
$ws.Range("H3").Value2 = 0.6963446507366129
$ws.Range("K3").Value2 = 0.2272977185796634
$ws.Range("M3").Value2 = 0.1926867908534149
$ws.Range("O3").Value2 = 2.49417829921704
$ws.Range("B4").Value2 = 0.2236698707330902
$ws.Range("C4").Value2 = 0.05527214612481401
$ws.Range("D4").Value2 = 0.06692929383598312
$ws.Range("E4").Value2 = 0.1510676827681579
$ws.Range("G4").Value2 = 0.5634908400214229
$ws.Range("H4").Value2 = 0.6990740303626097
$ws.Range("K4").Value2 = 0.2067641700347025
$ws.Range("M4").Value2 = 0.1791362881069887
$ws.Range("O4").Value2 = 2.50042740946779
$ws.Range("B5").Value2 = 0.2154280791049814
$ws.Range("C5").Value2 = 0.05448256808847418
$ws.Range("D5").Value2 = 0.06511441704545007
$ws.Range("E5").Value2 = 0.1484049154837024
$ws.Range("G5").Value2 = 0.5637290229708682
$ws.Range("H5").Value2 = 0.7002664422876848
$ws.Range("K5").Value2 = 0.1983806189188044
$ws.Range("M5").Value2 = 0.1736303727481143
$ws.Range("O5").Value2 = 2.503351329372492
$ws.Range("B6").Value2 = 0.2140594256084682
$ws.Range("C6").Value2 = 0.05435124168483441
$ws.Range("D6").Value2 = 0.06481356434686347
$ws.Range("E6").Value2 = 0.1479645583763229
$ws.Range("G6").Value2 = 0.5637746095106735
$ws.Range("H6").Value2 = 0.7004692838213629
$ws.Range("K6").Value2 = 0.1969875848294578
$ws.Range("M6").Value2 = 0.1727170866715682
$ws.Range("O6").Value2 = 2.503859627017704
$ws.Range("B7").Value2 = 0.2235587267143728
$ws.Range("C7").Value2 = 0.05526151221494047
$ws.Range("D7").Value2 = 0.06690478385058896
$ws.Range("E7").Value2 = 0.1510316513858427
$ws.Range("G7").Value2 = 0.5634936474712902
$ws.Range("H7").Value2 = 0.6990897870294006
$ws.Range("K7").Value2 = 0.2066511706122895
$ws.Range("M7").Value2 = 0.1790619685639214
$ws.Range("O7").Value2 = 2.500465314960337
$ws.Range("B8").Value2 = 0.2654581103816724
$ws.Range("C8").Value2 = 0.05924638379615033
$ws.Range("D8").Value2 = 0.07620632291454399
$ws.Range("E8").Value2 = 0.1648283296542132
$ws.Range("G8").Value2 = 0.5631323021978289
$ws.Range("H8").Value2 = 0.6937061287907369
$ws.Range("K8").Value2 = 0.2491531912021401
$ws.Range("M8").Value2 = 0.2072013069083951
$ws.Range("O8").Value2 = 2.488800616244646
$ws.Range("B9").Value2 = 0.3475310849370032
$ws.Range("C9").Value2 = 0.06694314097400422
$ws.Range("D9").Value2 = 0.09471028873775822
$ws.Range("E9").Value2 = 0.1928431101611139
$ws.Range("G9").Value2 = 0.5655920439862285
$ws.Range("H9").Value2 = 0.6856774116838125
$ws.Range("K9").Value2 = 0.3319683319628552
$ws.Range("M9").Value2 = 0.2628934251986763
$ws.Range("O9").Value2 = 2.477857979900108
$ws.Range("B10").Value2 = 0.4077610223201873
$ws.Range("C10").Value2 = 0.07252453541718751
$ws.Range("D10").Value2 = 0.1084674912474952
$ws.Range("E10").Value2 = 0.214029077322536
$ws.Range("G10").Value2 = 0.569342473338736
$ws.Range("H10").Value2 = 0.6813200748703281
$ws.Range("K10").Value2 = 0.3924754491919771
$ws.Range("M10").Value2 = 0.3041321384505693
$ws.Range("O10").Value2 = 2.477117403796029
$ws.Range("B11").Value2 = 0.4351437615016494
$ws.Range("C11").Value2 = 0.07504738402225541
$ws.Range("D11").Value2 = 0.1147616677293541
$ws.Range("E11").Value2 = 0.2238022839011649
$ws.Range("G11").Value2 = 0.5714731970388414
$ws.Range("H11").Value2 = 0.6796725134369836
$ws.Range("K11").Value2 = 0.4199260514892842
$ws.Range("M11").Value2 = 0.3229650030318254
$ws.Range("O11").Value2 = 2.478371022581854
$ws.Range("B12").Value2 = 0.4455102310436416
$ws.Range("C12").Value2 = 0.07600035889831247
$ws.Range("D12").Value2 = 0.1171502789751315
$ws.Range("E12").Value2 = 0.2275229270960679
$ws.Range("G12").Value2 = 0.5723412943671633
$ws.Range("H12").Value2 = 0.6790967385162787
$ws.Range("K12").Value2 = 0.4303098469221993
$ws.Range("M12").Value2 = 0.3301071429053977
$ws.Range("O12").Value2 = 2.479074843028712
$ws.Range("B13").Value2 = 0.4432777572003772
$ws.Range("C13").Value2 = 0.07579522504454417
$ws.Range("D13").Value2 = 0.1166356205002899
$ws.Range("E13").Value2 = 0.2267207381240581
$ws.Range("G13").Value2 = 0.5721516076737743
$ws.Range("H13").Value2 = 0.6792186016949984
$ws.Range("K13").Value2 = 0.4280740139880947
$ws.Range("M13").Value2 = 0.3285684860584936
$ws.Range("O13").Value2 = 2.478913067144191
$ws.Range("B14").Value2 = 0.4359966754980462
$ws.Range("C14").Value2 = 0.07512583368516346
$ws.Range("D14").Value2 = 0.114958077294375
$ws.Range("E14").Value2 = 0.224107986409102
$ws.Range("G14").Value2 = 0.5715433876092106
$ws.Range("H14").Value2 = 0.6796241796515119
$ws.Range("K14").Value2 = 0.4207805586873405
$ws.Range("M14").Value2 = 0.3235523797612103
$ws.Range("O14").Value2 = 2.478424332542914
$ws.Range("B15").Value2 = 0.4315364272563329
$ws.Range("C15").Value2 = 0.0747155016209291
$ws.Range("D15").Value2 = 0.1139312026224815
$ws.Range("E15").Value2 = 0.2225101787231409
$ws.Range("G15").Value2 = 0.5711788155439734
$ws.Range("H15").Value2 = 0.679878874733248
$ws.Range("K15").Value2 = 0.4163116409533529
$ws.Range("M15").Value2 = 0.3204812434373139
$ws.Range("O15").Value2 = 2.478154815870397
$ws.Range("B16").Value2 = 0.4059711317275685
$ws.Range("C16").Value2 = 0.07235933195445909
$ws.Range("D16").Value2 = 0.1080568732274116
$ws.Range("E16").Value2 = 0.2133931246938729
$ws.Range("G16").Value2 = 0.569211785381782
$ws.Range("H16").Value2 = 0.6814344804813999
$ws.Range("K16").Value2 = 0.3906799519434401
$ws.Range("M16").Value2 = 0.3029028405989536
$ws.Range("O16").Value2 = 2.47706751526303
$ws.Range("B17").Value2 = 0.3902831869671104
$ws.Range("C17").Value2 = 0.07090972461583078
$ws.Range("D17").Value2 = 0.1044623524562382
$ws.Range("E17").Value2 = 0.2078350111070932
$ws.Range("G17").Value2 = 0.5681139611133403
$ws.Range("H17").Value2 = 0.6824744993305814
$ws.Range("K17").Value2 = 0.3749363637939496
$ws.Range("M17").Value2 = 0.2921378148903884
$ws.Range("O17").Value2 = 2.476808136545941
$ws.Range("B18").Value2 = 0.3812584019926533
$ws.Range("C18").Value2 = 0.07007443112082967
$ws.Range("D18").Value2 = 0.1023982652047266
$ws.Range("E18").Value2 = 0.2046508803680638
$ws.Range("G18").Value2 = 0.5675224797436158
$ws.Range("H18").Value2 = 0.6831041857654014
$ws.Range("K18").Value2 = 0.3658740935098592
$ws.Range("M18").Value2 = 0.2859529423762268
$ws.Range("O18").Value2 = 2.476808631985847
$ws.Range("B19").Value2 = 0.3782025221682659
$ws.Range("C19").Value2 = 0.06979135596023411
$ws.Range("D19").Value2 = 0.1016999830986407
$ws.Range("E19").Value2 = 0.2035749707143992
$ws.Range("G19").Value2 = 0.5673290717667783
$ws.Range("H19").Value2 = 0.6833227959729697
$ws.Range("K19").Value2 = 0.3628045804751139
$ws.Range("M19").Value2 = 0.2838600327601526
$ws.Range("O19").Value2 = 2.476834496110001
$ws.Range("B20").Value2 = 0.3919533534053414
$ws.Range("C20").Value2 = 0.07106419523010743
$ws.Range("D20").Value2 = 0.1048446450838583
$ws.Range("E20").Value2 = 0.2084253605403319
$ws.Range("G20").Value2 = 0.5682266895826871
$ws.Range("H20").Value2 = 0.6823605278944171
$ws.Range("K20").Value2 = 0.3766130207957019
$ws.Range("M20").Value2 = 0.2932830574727703
$ws.Range("O20").Value2 = 2.476820253375308
$ws.Range("B21").Value2 = 0.4381353847454363
$ws.Range("C21").Value2 = 0.07532251495644005
$ws.Range("D21").Value2 = 0.1154506730330382
$ws.Range("E21").Value2 = 0.2248748771986371
$ws.Range("G21").Value2 = 0.5717203731875031
$ws.Range("H21").Value2 = 0.6795037455564312
$ws.Range("K21").Value2 = 0.4229231295781517
$ws.Range("M21").Value2 = 0.3250254454734502
$ws.Range("O21").Value2 = 2.478561664882733
$ws.Range("B22").Value2 = 0.4683016007395224
$ws.Range("C22").Value2 = 0.07809171772451862
$ws.Range("D22").Value2 = 0.122412290902119
$ws.Range("E22").Value2 = 0.235740771119751
$ws.Range("G22").Value2 = 0.5743607235752535
$ws.Range("H22").Value2 = 0.6779171601954772
$ws.Range("K22").Value2 = 0.4531242514794087
$ws.Range("M22").Value2 = 0.3458323924471003
$ws.Range("O22").Value2 = 2.481035340408681
$ws.Range("B23").Value2 = 0.4522029751958883
$ws.Range("C23").Value2 = 0.07661502624837624
$ws.Range("D23").Value2 = 0.1186940117456032
$ws.Range("E23").Value2 = 0.2299308196035525
$ws.Range("G23").Value2 = 0.5729187928541677
$ws.Range("H23").Value2 = 0.6787382847317929
$ws.Range("K23").Value2 = 0.4370114608734355
$ws.Range("M23").Value2 = 0.3347217024130984
$ws.Range("O23").Value2 = 2.479592761118852
$ws.Range("B24").Value2 = 0.3911982884957013
$ws.Range("C24").Value2 = 0.07099436496453393
$ws.Range("D24").Value2 = 0.1046718029525096
$ws.Range("E24").Value2 = 0.2081584283755831
$ws.Range("G24").Value2 = 0.5681756014764687
$ws.Range("H24").Value2 = 0.6824119554304247
$ws.Range("K24").Value2 = 0.3758550386194486
$ws.Range("M24").Value2 = 0.2927652805849448
$ws.Range("O24").Value2 = 2.476814309334571
$ws.Range("B25").Value2 = 0.3253393524736339
$ws.Range("C25").Value2 = 0.06487373216842229
$ws.Range("D25").Value2 = 0.0896760803344705
$ws.Range("E25").Value2 = 0.1851597280026809
$ws.Range("G25").Value2 = 0.5645862085312956
$ws.Range("H25").Value2 = 0.6875787031234495
$ws.Range("K25").Value2 = 0.3096228981082163
$ws.Range("M25").Value2 = 0.2477714741956163
$ws.Range("O25").Value2 = 2.479538280106482
